$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns(7).Insert()

$ws.Range("G3").Value = 50
$ws.Range("G4").Value = 50
$ws.Range("G4").Font.Bold = $true
$ws.Range("G4").Font.Italic = $true
$ws.Range("G5").Value = 50
$ws.Range("G7").Value = 50
$ws.Range("G8").Value = 50
$ws.Range("G9").Value = 50
$ws.Range("G12").Value = 90
$ws.Range("G13").Value = 90
$ws.Range("G14").Value = 90
$ws.Range("G16").Value = 90
$ws.Range("G17").Value = 90

$ws.Range("F1").Copy()
$ws.Range("G18").PasteSpecial(-4122)
